$d = $word.ActiveDocument
$sec = $d.Sections(1)

# ---------------------------------------------------------------------
# The document has "different first page" headers/footers turned on, so
# each of the two stories (header / footer) has two physical parts:
#   index 1 = default (odd/all other pages)
#   index 2 = first page
# Three inline pictures (Pearson logo x2, BTec logo x1) were renamed by
# the original edit:
#   footer (default)    Pearson logo : image1.png -> image2.png
#   footer (first page) Pearson logo : image1.png -> image2.png
#   header (first page) BTec logo    : image2.jpg -> image1.jpg
# ---------------------------------------------------------------------

$footerDefault = $sec.Footers(1).Range.InlineShapes(1)
$footerDefault.Name = "image2.png"

$footerFirst = $sec.Footers(2).Range.InlineShapes(1)
$footerFirst.Name = "image2.png"

$headerFirst = $sec.Headers(2).Range.InlineShapes(1)
$headerFirst.Name = "image1.jpg"
